# Implement 'remember me' and 'reset password' features:
# these TODO items are now done, so remove them (and the now-redundant
# "Consider to move to login page..." paragraph keeps its place as the
# surviving item) from the "High" priority To Do list.
#
# Before (High priority section):
#   ... Add feature to logically delete a user
#   User gets log out when he changes Profile (one possible solution is to
#     change `remember_token` only when email or password is changed)
#   Consider to move to login page instead of home page when access denied
#   Find out whether we can use `SecureRandom` module to generate random
#     password for newly created users
#   Implement "remember me" and "reset password" features
#   User should be redirected to the desired page after successful login ...
#
# After:
#   ... Add feature to logically delete a user
#   Consider to move to login page instead of home page when access denied
#   User should be redirected to the desired page after successful login ...

$d = $word.ActiveDocument

# Locate the three paragraphs to remove by matching distinctive (unique)
# substrings of their text. Each needle occurs in exactly one paragraph.
$targets = @(
    'remember me',
    'SecureRandom',
    'User gets log out when he changes Profile'
)

foreach ($needle in $targets) {
    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
        $para = $d.Paragraphs.Item($i)
        if ($para.Range.Text -match [regex]::Escape($needle)) {
            $para.Range.Delete()
            break
        }
    }
}
